$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.362.17"
$ws.Range("E2").Value = "  +8.81%  "
$ws.Range("D3").Value = "1.598.57"
$ws.Range("E3").Value = "  +8.03%  "
$ws.Range("E4").Value = "  -0.48%  "
$ws.Range("D5").Value = "'0.9945"
$ws.Range("D6").Value = "'288.06"
$ws.Range("E6").Value = "  +3.20%  "
$ws.Range("D7").Value = "'0.3684"
$ws.Range("E7").Value = "  +0.71%  "
$ws.Range("D8").Value = "'0.3385"
$ws.Range("E8").Value = "  +10.04%  "
$ws.Range("D9").Value = "'42.68"
$ws.Range("E9").Value = "  +6.62%  "
$ws.Range("D10").Value = "'1.138"
$ws.Range("E10").Value = "  +7.53%  "
$ws.Range("D11").Value = "'0.07031"
$ws.Range("E11").Value = "  +5.57%  "
$ws.Range("D12").Value = "'1.000"
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("D13").Value = "'19.69"
$ws.Range("E13").Value = "  +9.02%  "
$ws.Range("D14").Value = "'5.916"
$ws.Range("E14").Value = "  +7.42%  "
$ws.Range("D15").Value = "'6.617"
$ws.Range("E15").Value = "  +6.67%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.00001082"
$ws.Range("E16").Value = "  +5.31%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "1.596.77"
$ws.Range("E17").Value = "  +7.76%  "
$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").Value = "'0.9937"
$ws.Range("E18").Value = "  +1.64%  "
$ws.Range("D19").Value = "'0.06594"
$ws.Range("E19").Value = "  +11.02%  "
$ws.Range("D20").Value = "'77.93"
$ws.Range("E20").Value = "  +12.26%  "
$ws.Range("E21").Value = "  +11.17%  "
$ws.Range("D22").Value = "'6.009"
$ws.Range("E22").Value = "  +9.51%  "
$ws.Range("D23").Value = "'11.76"
$ws.Range("E23").Value = "  +6.60%  "
$ws.Range("D24").Value = "22.366.40"
$ws.Range("E24").Value = "  +8.49%  "
$ws.Range("D25").Value = "'2.404"
$ws.Range("E25").Value = "  +6.44%  "
$ws.Range("D26").Value = "'2.507"
$ws.Range("E26").Value = "  +16.40%  "
$ws.Range("D27").Value = "'149.48"
$ws.Range("E27").Value = "  +5.66%  "
$ws.Range("D28").Value = "'19.57"
$ws.Range("E28").Value = "  +13.30%  "
$ws.Range("D29").Value = "1.773.59"
$ws.Range("E29").Value = "  +8.06%  "
$ws.Range("D30").Value = "'120.12"
$ws.Range("E30").Value = "  +5.60%  "
$ws.Range("D31").Value = "'4.208"
$ws.Range("E31").Value = "  +6.75%  "
$ws.Range("D32").Value = "'6.004"
$ws.Range("E32").Value = "  +19.66%  "
$ws.Range("D33").Value = "'0.9462"
$ws.Range("E33").Value = "  +15.64%  "
$ws.Range("D34").Value = "'0.08250"
$ws.Range("E34").Value = "  +3.03%  "
$ws.Range("D35").Value = "'1.612"
$ws.Range("E35").Value = "  +5.59%  "
$ws.Range("D36").Value = "'5.285"
$ws.Range("E36").Value = "  +11.73%  "
$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D37").Value = "'8.648"
$ws.Range("E37").Value = "  +12.04%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").Value = "'11.78"
$ws.Range("E38").Value = "  +12.70%  "
$ws.Range("D39").Value = "'0.06112"
$ws.Range("E39").Value = "  +5.33%  "
$ws.Range("D40").Value = "'1.240"
$ws.Range("E40").Value = "  +0.70%  "
$ws.Range("D41").Value = "'0.02214"
$ws.Range("E41").Value = "  +8.39%  "
$ws.Range("D42").Value = "'0.2021"
$ws.Range("E42").Value = "  +6.97%  "
$ws.Range("D43").Value = "'0.9937"
$ws.Range("E43").Value = "  +1.69%  "
$ws.Range("D44").Value = "'0.5889"
$ws.Range("E44").Value = "  +11.23%  "
$ws.Range("D45").Value = "'13.19"
$ws.Range("E45").Value = "  +8.51%  "
$ws.Range("D46").Value = "'3.665"
$ws.Range("E46").Value = "  +3.89%  "
$ws.Range("D47").Value = "'0.5688"
$ws.Range("E47").Value = "  +9.39%  "
$ws.Range("D48").Value = "'126.17"
$ws.Range("E48").Value = "  +5.78%  "
$ws.Range("E49").Value = "  +8.89%  "
$ws.Range("D50").Value = "'0.06802"
$ws.Range("E50").Value = "  +5.00%  "
$ws.Range("D51").Value = "'73.53"
$ws.Range("E51").Value = "  +8.72%  "
